$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (old column B) entirely, shifting
# all the columns to its right one place to the left.
$ws.Range("B1").EntireColumn.Delete()

# Append ".deja.deja.deja" to each of the (now shifted) header labels,
# i.e. what used to be C1:I1 and is now B1:H1.
$ws.Range("B1").Value = $ws.Range("B1").Value2 + ".deja.deja.deja"
$ws.Range("C1").Value = $ws.Range("C1").Value2 + ".deja.deja.deja"
$ws.Range("D1").Value = $ws.Range("D1").Value2 + ".deja.deja.deja"
$ws.Range("E1").Value = $ws.Range("E1").Value2 + ".deja.deja.deja"
$ws.Range("F1").Value = $ws.Range("F1").Value2 + ".deja.deja.deja"
$ws.Range("G1").Value = $ws.Range("G1").Value2 + ".deja.deja.deja"
$ws.Range("H1").Value = $ws.Range("H1").Value2 + ".deja.deja.deja"
